# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-26 07:15:21
#
# Re-applies the refreshed "Recorded By" ordering, the re-computed Missing /
# Pending session counters, and the status flip for the Pathology
# LAB/MUSEUM session (row 22) that moved from "Pending" to "Not Recorded".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" (column G) lists re-ordered by the refreshed sync ---
$ws.Range("G2").Value  = "Veronia.rafat@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value  = "Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value  = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G5").Value  = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G7").Value  = "Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G9").Value  = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- Recomputed summary counters (Missing Sessions / Pending Sessions) ---
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 14
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 14

# --- Row 22 (PATHOLOGY LAB/MUSEUM, session 1) flips from Pending to Not
#     Recorded now that its date has passed without a record. Re-style
#     A22:I22 to the "Not Recorded" look already used on row 29 (pink),
#     matching the refreshed status-color convention. ---
$ws.Range("A29:I29").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I22").Value = "Not Recorded"
